# Update the "Corr/total marks" on the concise marksheet:
#  - B11 (Marking -> Right column): 3 -> 5
#  - B12 (Total -> Right column): 60 -> 100
#  - E12 (Total -> Max column, displayed score string): "60/84" -> "100/140"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 100
$ws.Range("E12").Value = "100/140"
